# ============================================================================
# Adds "Player Info" and "ODI Batting Extra" sheets, and updates the
# "ODI Batting" sheet (rename MATCH_CARD_LINK -> MATCH_CODE, replace the
# full scorecard URL with just the numeric match code, drop the stray
# empty B11 cell).
#
# Row data below is encoded with "|" as the field separator and "~" as a
# sentinel for a genuinely empty field, to stay robust against any
# whitespace collapsing of the script source.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 0. Identify / normalise the pre-existing sheet.
# ---------------------------------------------------------------------------
$odiBatting = $wb.ActiveSheet
$odiBatting.Name = "ODI Batting"

# ---------------------------------------------------------------------------
# 1. "Player Info" sheet -- inserted BEFORE "ODI Batting".
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$playerInfoData = @"
ID|NAME|BATTING_HAND|BOWL_STYLE
3842|Usman T Khawaja|Left Handed|Right Arm Medium
"@

$piLines = $playerInfoData -split "`n"

# Header row (row 1), bold/bordered/centered like the rest of the workbook.
$piHeader = $piLines[0] -split "\|"
for ($c = 0; $c -lt $piHeader.Count; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $piHeader[$c].Trim()
    $cell.NumberFormat = "General"
}
$piHeaderRange = $playerInfo.Range("A1:D1")
$piHeaderRange.Font.Bold = $true
$piHeaderRange.HorizontalAlignment = -4108
$piHeaderRange.VerticalAlignment = -4160
$piHeaderRange.Borders.LineStyle = 1

# Data rows.
for ($r = 1; $r -lt $piLines.Count; $r++) {
    $line = $piLines[$r].Trim()
    if ($line -eq "") { continue }
    $vals = $line -split "\|"
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $cell = $playerInfo.Cells.Item($r + 1, $c + 1)
        $cell.NumberFormat = "@"
        $v = $vals[$c].Trim()
        if ($v -eq "~") { $v = "" }
        $cell.Value = $v
        $cell.NumberFormat = "General"
    }
}

$playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. "ODI Batting Extra" sheet -- inserted AFTER "ODI Batting".
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiExtra = $wb.Worksheets.Add($null, $odiBatting)
$odiExtra.Name = "ODI Batting Extra"

$odiExtraData = @"
MATCH_CODE|BATTING_POSITION|NUM_4|NUM_6|PERCENT_RUNS_OF_TOTAL|MAN_OF_MATCH
4236|3|2|0|14.78%|NO
4258|~|~|~|~|NO
4263|2|6|0|15.70%|NO
4266|2|11|1|33.23%|YES
4268|2|7|0|25.35%|NO
4270|1|10|2|36.76%|YES
4273|1|3|0|8.54%|NO
4274|1|8|0|30.88%|NO
4275|1|0|0|~|NO
4276|~|~|~|~|NO
4277|1|10|0|29.97%|NO
4306|3|1|0|7.18%|NO
4312|3|2|0|4.51%|NO
4316|4|4|1|13.29%|NO
4319|6|3|0|5.86%|NO
4322|3|1|0|2.99%|NO
4329|~|~|~|~|NO
4336|3|1|0|8.07%|NO
4341|3|5|0|36.21%|NO
4351|3|3|0|5.71%|NO
"@

$oeLines = $odiExtraData -split "`n"

# Header row (row 1).
$oeHeader = $oeLines[0] -split "\|"
for ($c = 0; $c -lt $oeHeader.Count; $c++) {
    $cell = $odiExtra.Cells.Item(1, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $oeHeader[$c].Trim()
    $cell.NumberFormat = "General"
}
$oeHeaderRange = $odiExtra.Range("A1:F1")
$oeHeaderRange.Font.Bold = $true
$oeHeaderRange.HorizontalAlignment = -4108
$oeHeaderRange.VerticalAlignment = -4160
$oeHeaderRange.Borders.LineStyle = 1

# Data rows. Column B (BATTING_POSITION) is numeric when present; the rest
# of the columns are kept as plain text, matching the source data.
for ($r = 1; $r -lt $oeLines.Count; $r++) {
    $line = $oeLines[$r].Trim()
    if ($line -eq "") { continue }
    $vals = $line -split "\|"
    if ($vals.Count -lt 6) { continue }
    $rowNum = $r + 1

    $cellA = $odiExtra.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $vals[0].Trim()
    $cellA.NumberFormat = "General"

    $bVal = $vals[1].Trim()
    $cellB = $odiExtra.Cells.Item($rowNum, 2)
    if ($bVal -eq "~") {
        $cellB.Value = ""
    } else {
        $cellB.NumberFormat = "General"
        $cellB.Value = [double]$bVal
    }

    for ($c = 2; $c -lt 6; $c++) {
        $cell = $odiExtra.Cells.Item($rowNum, $c + 1)
        $cell.NumberFormat = "@"
        $v = $vals[$c].Trim()
        if ($v -eq "~") { $v = "" }
        $cell.Value = $v
        $cell.NumberFormat = "General"
    }
}

$odiExtra.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. Update "ODI Batting": rename the MATCH_CARD_LINK header, replace the
#    full URL in column D with just the numeric match code, and drop the
#    stray empty B11 cell.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$d1 = $odiBatting.Cells.Item(1, 4)
$d1.NumberFormat = "@"
$d1.Value = "MATCH_CODE"
$d1.NumberFormat = "General"

for ($r = 2; $r -le 41; $r++) {
    $cell = $odiBatting.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
        $cell.NumberFormat = "General"
    }
}

$odiBatting.Cells.Item(11, 2).ClearContents()

$odiBatting.Range("A1").Select()

Write-Host "Workbook restructure complete."
